$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 29.51793756855773
$ws.Cells.Item(2, 3).Value = 8.231041776144481
$ws.Cells.Item(2, 4).Value = 12.73262989164373
$ws.Cells.Item(2, 5).Value = 11.75721815040084
$ws.Cells.Item(2, 7).Value = 3.902069260396217
$ws.Cells.Item(2, 10).Value = 7.402687214562372
$ws.Cells.Item(2, 11).Value = 24.75949731911664
$ws.Cells.Item(2, 12).Value = 13.70527809335086
$ws.Cells.Item(2, 13).Value = 24.03335456071629
$ws.Cells.Item(2, 14).Value = 32.93702529482636
$ws.Cells.Item(3, 2).Value = 29.45590054447274
$ws.Cells.Item(3, 3).Value = 8.138625573431337
$ws.Cells.Item(3, 4).Value = 12.74395807971209
$ws.Cells.Item(3, 5).Value = 11.77841584113204
$ws.Cells.Item(3, 7).Value = 3.906131233762383
$ws.Cells.Item(3, 10).Value = 7.397150935697447
$ws.Cells.Item(3, 11).Value = 24.72093980398493
$ws.Cells.Item(3, 12).Value = 13.72880034044123
$ws.Cells.Item(3, 13).Value = 24.05348281048499
$ws.Cells.Item(3, 14).Value = 32.87116483055614
$ws.Cells.Item(4, 2).Value = 29.42421008718831
$ws.Cells.Item(4, 3).Value = 8.0838317434093
$ws.Cells.Item(4, 4).Value = 12.75297283732156
$ws.Cells.Item(4, 5).Value = 11.79238910498318
$ws.Cells.Item(4, 7).Value = 3.908754214424401
$ws.Cells.Item(4, 10).Value = 7.393678037850587
$ws.Cells.Item(4, 11).Value = 24.70249044663355
$ws.Cells.Item(4, 12).Value = 13.74500692672137
$ws.Cells.Item(4, 13).Value = 24.06980769068407
$ws.Cells.Item(4, 14).Value = 32.83165419688608
$ws.Cells.Item(5, 2).Value = 29.41291251248413
$ws.Cells.Item(5, 3).Value = 8.062017456161277
$ws.Cells.Item(5, 4).Value = 12.75716415609244
$ws.Cells.Item(5, 5).Value = 11.79832471092138
$ws.Cells.Item(5, 7).Value = 3.909855641487255
$ws.Cells.Item(5, 10).Value = 7.392244179716107
$ws.Cells.Item(5, 11).Value = 24.69629002569542
$ws.Cells.Item(5, 12).Value = 13.75205491075088
$ws.Cells.Item(5, 13).Value = 24.07745730691392
$ws.Cells.Item(5, 14).Value = 32.81578923611274
$ws.Cells.Item(6, 2).Value = 29.41113437282712
$ws.Cells.Item(6, 3).Value = 8.058426979624466
$ws.Cells.Item(6, 4).Value = 12.75789138674798
$ws.Cells.Item(6, 5).Value = 11.79932490893149
$ws.Cells.Item(6, 7).Value = 3.910040501687837
$ws.Cells.Item(6, 10).Value = 7.392004960395743
$ws.Cells.Item(6, 11).Value = 24.69534014034782
$ws.Cells.Item(6, 12).Value = 13.75325202323008
$ws.Cells.Item(6, 13).Value = 24.07878773138442
$ws.Cells.Item(6, 14).Value = 32.81316919475224
$ws.Cells.Item(7, 2).Value = 29.42405117061537
$ws.Cells.Item(7, 3).Value = 8.083535433341062
$ws.Cells.Item(7, 4).Value = 12.75302726681326
$ws.Cells.Item(7, 5).Value = 11.79246817654916
$ws.Cells.Item(7, 7).Value = 3.908768936732994
$ws.Cells.Item(7, 10).Value = 7.393658775897542
$ws.Cells.Item(7, 11).Value = 24.70240148530287
$ws.Cells.Item(7, 12).Value = 13.74510018174583
$ws.Cells.Item(7, 13).Value = 24.06990681937393
$ws.Cells.Item(7, 14).Value = 32.83143927795227
$ws.Cells.Item(8, 2).Value = 29.49522343920398
$ws.Cells.Item(8, 3).Value = 8.198786842414421
$ws.Cells.Item(8, 4).Value = 12.73610845076808
$ws.Cells.Item(8, 5).Value = 11.76432864894865
$ws.Cells.Item(8, 7).Value = 3.903443148419076
$ws.Cells.Item(8, 10).Value = 7.400793510215396
$ws.Cells.Item(8, 11).Value = 24.74512087538977
$ws.Cells.Item(8, 12).Value = 13.71302256996461
$ws.Cells.Item(8, 13).Value = 24.03947156935403
$ws.Cells.Item(8, 14).Value = 32.9141255993999
$ws.Cells.Item(9, 2).Value = 29.68520641213402
$ws.Cells.Item(9, 3).Value = 8.439148686877104
$ws.Cells.Item(9, 4).Value = 12.71926973327787
$ws.Cells.Item(9, 5).Value = 11.71672271614159
$ws.Cells.Item(9, 7).Value = 3.894016454912447
$ws.Cells.Item(9, 10).Value = 7.414206369540855
$ws.Cells.Item(9, 11).Value = 24.87013809102669
$ws.Cells.Item(9, 12).Value = 13.66410756629502
$ws.Cells.Item(9, 13).Value = 24.01125953358175
$ws.Cells.Item(9, 14).Value = 33.08352211726717
$ws.Cells.Item(10, 2).Value = 29.85493800452128
$ws.Cells.Item(10, 3).Value = 8.622937477400679
$ws.Cells.Item(10, 4).Value = 12.71685474121581
$ws.Cells.Item(10, 5).Value = 11.68633212282512
$ws.Cells.Item(10, 7).Value = 3.887702771987468
$ws.Cells.Item(10, 10).Value = 7.423716614174706
$ws.Cells.Item(10, 11).Value = 24.98678562526958
$ws.Cells.Item(10, 12).Value = 13.63668955107874
$ws.Cells.Item(10, 13).Value = 24.00971185055052
$ws.Cells.Item(10, 14).Value = 33.2122444769352
$ws.Cells.Item(11, 2).Value = 29.93854729575154
$ws.Cells.Item(11, 3).Value = 8.707768910614469
$ws.Cells.Item(11, 4).Value = 12.71791467382224
$ws.Cells.Item(11, 5).Value = 11.67349544317535
$ws.Cells.Item(11, 7).Value = 3.884961722510834
$ws.Cells.Item(11, 10).Value = 7.427970157557922
$ws.Cells.Item(11, 11).Value = 25.04514153292727
$ws.Cells.Item(11, 12).Value = 13.62606419966451
$ws.Cells.Item(11, 13).Value = 24.01316613049396
$ws.Cells.Item(11, 14).Value = 33.27169499452891
$ws.Cells.Item(12, 2).Value = 29.97111272321605
$ws.Cells.Item(12, 3).Value = 8.740039698826653
$ws.Cells.Item(12, 4).Value = 12.7186258812754
$ws.Cells.Item(12, 5).Value = 11.66877608470082
$ws.Cells.Item(12, 7).Value = 3.883942475885269
$ws.Cells.Item(12, 10).Value = 7.429570538440665
$ws.Cells.Item(12, 11).Value = 25.06799073163857
$ws.Cells.Item(12, 12).Value = 13.62230604755708
$ws.Cells.Item(12, 13).Value = 24.01507105501909
$ws.Cells.Item(12, 14).Value = 33.29433247152322
$ws.Cells.Item(13, 2).Value = 29.96405921355767
$ws.Cells.Item(13, 3).Value = 8.733083534017066
$ws.Cells.Item(13, 4).Value = 12.71845894044894
$ws.Cells.Item(13, 5).Value = 11.66978619231949
$ws.Cells.Item(13, 7).Value = 3.884161157785819
$ws.Cells.Item(13, 10).Value = 7.429226326032684
$ws.Cells.Item(13, 11).Value = 25.06303651196631
$ws.Cells.Item(13, 12).Value = 13.62310363083785
$ws.Cells.Item(13, 13).Value = 24.01463427024143
$ws.Cells.Item(13, 14).Value = 33.28945159342575
$ws.Cells.Item(14, 2).Value = 29.94120844117557
$ws.Cells.Item(14, 3).Value = 8.71042106905276
$ws.Cells.Item(14, 4).Value = 12.71796697944549
$ws.Cells.Item(14, 5).Value = 11.67310434331256
$ws.Cells.Item(14, 7).Value = 3.884877493795622
$ws.Cells.Item(14, 10).Value = 7.428102028514734
$ws.Cells.Item(14, 11).Value = 25.04700635787474
$ws.Cells.Item(14, 12).Value = 13.62574969567505
$ws.Cells.Item(14, 13).Value = 24.01331089503813
$ws.Cells.Item(14, 14).Value = 33.27355491286701
$ws.Cells.Item(15, 2).Value = 29.92732897952274
$ws.Cells.Item(15, 3).Value = 8.696557930173899
$ws.Cells.Item(15, 4).Value = 12.71770596970971
$ws.Cells.Item(15, 5).Value = 11.6751552346356
$ws.Cells.Item(15, 7).Value = 3.885318705862289
$ws.Cells.Item(15, 10).Value = 7.42741202010579
$ws.Cells.Item(15, 11).Value = 25.03728492692688
$ws.Cells.Item(15, 12).Value = 13.62740504837708
$ws.Cells.Item(15, 13).Value = 24.01257798137529
$ws.Cells.Item(15, 14).Value = 33.26383388431336
$ws.Cells.Item(16, 2).Value = 29.84960137569577
$ws.Cells.Item(16, 3).Value = 8.617415696381736
$ws.Cells.Item(16, 4).Value = 12.71682886249793
$ws.Cells.Item(16, 5).Value = 11.68719087293277
$ws.Cells.Item(16, 7).Value = 3.887884533401446
$ws.Cells.Item(16, 10).Value = 7.423437181540008
$ws.Cells.Item(16, 11).Value = 24.98307759676472
$ws.Cells.Item(16, 12).Value = 13.63742109818791
$ws.Cells.Item(16, 13).Value = 24.0095696938981
$ws.Cells.Item(16, 14).Value = 33.20837692305491
$ws.Cells.Item(17, 2).Value = 29.8035459210821
$ws.Cells.Item(17, 3).Value = 8.569157457613494
$ws.Cells.Item(17, 4).Value = 12.71684328128535
$ws.Cells.Item(17, 5).Value = 11.69482709297658
$ws.Cells.Item(17, 7).Value = 3.889492072287444
$ws.Cells.Item(17, 10).Value = 7.420980203092455
$ws.Cells.Item(17, 11).Value = 24.95117188498257
$ws.Cells.Item(17, 12).Value = 13.64403861538351
$ws.Cells.Item(17, 13).Value = 24.00878856167249
$ws.Cells.Item(17, 14).Value = 33.17458227453182
$ws.Cells.Item(18, 2).Value = 29.7776592242753
$ws.Cells.Item(18, 3).Value = 8.541517644273863
$ws.Cells.Item(18, 4).Value = 12.71705478734851
$ws.Cells.Item(18, 5).Value = 11.69931228305198
$ws.Cells.Item(18, 7).Value = 3.890429030790445
$ws.Cells.Item(18, 10).Value = 7.419560163121131
$ws.Cells.Item(18, 11).Value = 24.93331936389889
$ws.Cells.Item(18, 12).Value = 13.64801872086518
$ws.Cells.Item(18, 13).Value = 24.00873076024753
$ws.Cells.Item(18, 14).Value = 33.15522887234898
$ws.Cells.Item(19, 2).Value = 29.76899848312876
$ws.Cells.Item(19, 3).Value = 8.532180272482718
$ws.Cells.Item(19, 4).Value = 12.71716131642724
$ws.Cells.Item(19, 5).Value = 11.70084688549613
$ws.Cells.Item(19, 7).Value = 3.890748392688532
$ws.Cells.Item(19, 10).Value = 7.419078181894919
$ws.Cells.Item(19, 11).Value = 24.92736076366198
$ws.Cells.Item(19, 12).Value = 13.64939618640786
$ws.Cells.Item(19, 13).Value = 24.00877845647855
$ws.Cells.Item(19, 14).Value = 33.14869074678661
$ws.Cells.Item(20, 2).Value = 29.80838628429242
$ws.Cells.Item(20, 3).Value = 8.574282727145478
$ws.Cells.Item(20, 4).Value = 12.71682071825049
$ws.Cells.Item(20, 5).Value = 11.69400457871362
$ws.Cells.Item(20, 7).Value = 3.889319670270255
$ws.Cells.Item(20, 10).Value = 7.421242460873088
$ws.Cells.Item(20, 11).Value = 24.95451674753147
$ws.Cells.Item(20, 12).Value = 13.64331617429021
$ws.Cells.Item(20, 13).Value = 24.00883120438381
$ws.Cells.Item(20, 14).Value = 33.17817106535168
$ws.Cells.Item(21, 2).Value = 29.94789585206126
$ws.Cells.Item(21, 3).Value = 8.71707382766472
$ws.Cells.Item(21, 4).Value = 12.71810307665178
$ws.Cells.Item(21, 5).Value = 11.6721258824025
$ws.Cells.Item(21, 7).Value = 3.884666581185016
$ws.Cells.Item(21, 10).Value = 7.428432541778399
$ws.Cells.Item(21, 11).Value = 25.05169450262919
$ws.Cells.Item(21, 12).Value = 13.6249652800751
$ws.Cells.Item(21, 13).Value = 24.01368341385232
$ws.Cells.Item(21, 14).Value = 33.27822079502418
$ws.Cells.Item(22, 2).Value = 30.04433685948529
$ws.Cells.Item(22, 3).Value = 8.811239296663114
$ws.Cells.Item(22, 4).Value = 12.72074676471895
$ws.Cells.Item(22, 5).Value = 11.65865210080914
$ws.Cells.Item(22, 7).Value = 3.881734632251082
$ws.Cells.Item(22, 10).Value = 7.433071436439866
$ws.Cells.Item(22, 11).Value = 25.11957837901032
$ws.Cells.Item(22, 12).Value = 13.61451897336424
$ws.Cells.Item(22, 13).Value = 24.02033288327882
$ws.Cells.Item(22, 14).Value = 33.34433500958566
$ws.Cells.Item(23, 2).Value = 29.99238834364798
$ws.Cells.Item(23, 3).Value = 8.760913854737435
$ws.Cells.Item(23, 4).Value = 12.71917078713579
$ws.Cells.Item(23, 5).Value = 11.66576796287252
$ws.Cells.Item(23, 7).Value = 3.883289523966084
$ws.Cells.Item(23, 10).Value = 7.430601045796315
$ws.Cells.Item(23, 11).Value = 25.08295091932898
$ws.Cells.Item(23, 12).Value = 13.61995287923978
$ws.Cells.Item(23, 13).Value = 24.01646610991732
$ws.Cells.Item(23, 14).Value = 33.30898344996633
$ws.Cells.Item(24, 2).Value = 29.80619611347262
$ws.Cells.Item(24, 3).Value = 8.571965265109242
$ws.Cells.Item(24, 4).Value = 12.71683028598336
$ws.Cells.Item(24, 5).Value = 11.69437614139183
$ws.Cells.Item(24, 7).Value = 3.889397573459616
$ws.Cells.Item(24, 10).Value = 7.42112391746948
$ws.Cells.Item(24, 11).Value = 24.95300300663943
$ws.Cells.Item(24, 12).Value = 13.64364224291835
$ws.Cells.Item(24, 13).Value = 24.00881070677791
$ws.Cells.Item(24, 14).Value = 33.1765483372916
$ws.Cells.Item(25, 2).Value = 29.62846831051642
$ws.Cells.Item(25, 3).Value = 8.372745962413619
$ws.Cells.Item(25, 4).Value = 12.72207578115998
$ws.Cells.Item(25, 5).Value = 11.72879376823655
$ws.Cells.Item(25, 7).Value = 3.896458557863818
$ws.Cells.Item(25, 10).Value = 7.410638881777902
$ws.Cells.Item(25, 11).Value = 24.83193379999848
$ws.Cells.Item(25, 12).Value = 13.67584349053365
$ws.Cells.Item(25, 13).Value = 24.01552199728835
$ws.Cells.Item(25, 14).Value = 33.03693492227386
